# 2661.firstCompleteIndex & 2662.minimumCost & lcp.79.extractMantra
# Adds a new worksheet "2662. 前往目标的最小代价" (after the existing
# "2106. 摘水果" sheet) containing the minimum-cost-path trace table,
# and leaves it as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet right after sheet1; Excel auto-activates it,
#     which also clears tabSelected on sheet1 and sets activeTab on the
#     workbook view - matching the diff automatically.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2662. 前往目标的最小代价"

# --- row 1: merged title banner (mirrors the A1:N1 banner on sheet1) -
#     reuse sheet1's banner format (center + wrap) via copy/paste-formats
#     so we don't fork a brand-new duplicate style entry.
$null = $ws1.Range("A1").Copy()
$null = $ws2.Range("A1:Q1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A1").Value = "[1,1] [10,4] [[4,2,1,1,3],[1,2,7,4,4],[10,3,6,1,2],[6,1,1,2,3]]"
$ws2.Range("A1:Q1").Merge()

# --- row 3: column headers for the trace grid
$ws2.Range("B3").Value = "1,1"
$ws2.Range("C3").Value = "10,4"
$ws2.Range("D3").Value = "4,2"
$ws2.Range("E3").Value = "1,2"
$ws2.Range("F3").Value = "7,4"
$ws2.Range("G3").Value = "10,3"
$ws2.Range("H3").Value = "6,1"

# --- row 4
$ws2.Range("A4").Value = 1.1
$ws2.Range("B4").Value = 0
$ws2.Range("B4").Font.Color = 255   # red - on the traced path
$ws2.Range("C4").Value = 12
$ws2.Range("D4").Value = 4
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 9
$ws2.Range("G4").Value = 11
$ws2.Range("H4").Value = 5

# --- row 5
$ws2.Range("A5").Value = 1.2
$ws2.Range("C5").Value = 12
$ws2.Range("D5").Value = 3
$ws2.Range("E5").Value = 1
$ws2.Range("E5").Font.Color = 255   # red - on the traced path
$ws2.Range("F5").Value = 9
$ws2.Range("G5").Value = 11
$ws2.Range("H5").Value = 5

# --- row 6
$ws2.Range("A6").Value = 4.2
$ws2.Range("C6").Value = 12
$ws2.Range("D6").Value = 3
$ws2.Range("D6").Font.Color = 255   # red - on the traced path
$ws2.Range("F6").Value = 5
$ws2.Range("G6").Value = 11
$ws2.Range("H6").Value = 4

# --- row 7
$ws2.Range("A7").Value = 6.1
$ws2.Range("C7").Value = 11
$ws2.Range("F7").Value = 5
$ws2.Range("G7").Value = 11
$ws2.Range("H7").Value = 4
$ws2.Range("H7").Font.Color = 255   # red - on the traced path

# --- row 8
$ws2.Range("C8").Value = 11
$ws2.Range("F8").Value = 5
$ws2.Range("G8").Value = 6

# --- rows 10-14: legend mapping each step back to its cost
$ws2.Range("A10").Value = "1,1"
$ws2.Range("B10").Value = 0
$ws2.Range("A11").Value = "1,2"
$ws2.Range("B11").Value = 1
$ws2.Range("A12").Value = "6,1"
$ws2.Range("B12").Value = 3
$ws2.Range("A13").Value = "10,3"
$ws2.Range("B13").Value = 2
$ws2.Range("A14").Value = "10,4"
$ws2.Range("B14").Value = 1

$ws2.Range("O27").Select()

Write-Host "2662 sheet created"
